$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cell content (sharedStrings additions, indices 11-17) ---
$ws.Range("C6").Value = "encrypt('66' || SERVICENUMBER)"
$ws.Range("B7").Value = "CBS_TS_PRODPRICEPLN_DTN"
$ws.Range("C7").Value = "DIM_CBS _TS_PRODPRICEPLN"
$ws.Range("B8").Value = "CBS_PRODUCTEXTATTR_DTN"
$ws.Range("C8").Value = "DIM_CBS_PRODUCTEXTATTR"
$ws.Range("B9").Value = "CBS_USAGESERVICETPE_DTN"
$ws.Range("C9").Value = "DIM_DEF_CBS_USAGESERVICETYPE"

# --- Fonts / styles ---
$ws.Range("B7:B8").Font.Name = "Calibri"
$ws.Range("B7:B8").Font.Size = 8
$ws.Range("B7:B8").Font.Color = 255

$ws.Range("C7:C8").Font.Name = "Calibri Light"
$ws.Range("C7:C8").Font.Size = 12
$ws.Range("C7:C8").Font.Color = 4486484
$ws.Range("C7:C8").HorizontalAlignment = -4131
$ws.Range("C7:C8").VerticalAlignment = -4108

$ws.Range("C9").Font.Name = "Calibri"
$ws.Range("C9").Font.Size = 12
$ws.Range("C9").Font.Color = 4486484
$ws.Range("C9").HorizontalAlignment = -4131
$ws.Range("C9").VerticalAlignment = -4108
